$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = -2
    4  = -1
    5  = -3
    6  = -1
    7  = -2
    8  = -1
    10 = -7
    11 = -2
    12 = 3
    13 = -1
    14 = -2
    15 = 1
    16 = -5
    17 = 6
    18 = -5
    19 = -3
    21 = -5
    22 = 5
    24 = 7
    25 = -2
    26 = 1
    27 = 4
    28 = 7
    29 = 10
    31 = 2
    32 = 7
    33 = 8
    34 = -2
    36 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
